{"js": "// Replace the arithmetic-problem text in each table cell with the new\n// problem/answer pair. Every \"before\" string occurs exactly once in the\n// document, so a direct search + replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"355\u00d78=2840\", \"562\u00d76=3372\"],\n  [\"673\u00d77=4711\", \"139\u00d76=834\"],\n  [\"573\u00d76=3438\", \"872\u00d79=7848\"],\n  [\"230\u00d72=460\", \"423\u00d76=2538\"],\n  [\"638\u00d78=5104\", \"739\u00d75=3695\"],\n  [\"216\u00d74=864\", \"880\u00d76=5280\"],\n  [\"599\u00d73=1797\", \"183\u00d74=732\"],\n  [\"824\u00d72=1648\", \"294\u00d79=2646\"],\n  [\"787\u00d78=6296\", \"755\u00d76=4530\"],\n  [\"815\u00d74=3260\", \"837\u00d72=1674\"],\n  [\"398\u00d72=796\", \"663\u00d74=2652\"],\n  [\"271\u00d77=1897\", \"634\u00d78=5072\"],\n  [\"943\u00d79=8487\", \"611\u00d73=1833\"],\n  [\"762\u00d76=4572\", \"828\u00d74=3312\"],\n  [\"406\u00d77=2842\", \"958\u00d72=1916\"],\n  [\"494\u00d78=3952\", \"296\u00d74=1184\"],\n  [\"924\u00d79=8316\", \"165\u00d78=1320\"],\n  [\"307\u00d79=2763\", \"719\u00d75=3595\"],\n  [\"239\u00d72=478\", \"551\u00d79=4959\"],\n  [\"846\u00d72=1692\", \"664\u00d77=4648\"],\n  [\"589\u00d72=1178\", \"347\u00d78=2776\"],\n  [\"740\u00d79=6660\", \"649\u00d77=4543\"],\n  [\"460\u00d74=1840\", \"569\u00d75=2845\"],\n  [\"439\u00d74=1756\", \"240\u00d73=720\"],\n  [\"391\u00d72=782\", \"313\u00d72=626\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the arithmetic-problem text in each table cell with the new\n# problem/answer pair. Every \"before\" string occurs exactly once in the\n# document, so Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"355\u00d78=2840\";  After = \"562\u00d76=3372\" },\n    @{ Before = \"673\u00d77=4711\";  After = \"139\u00d76=834\" },\n    @{ Before = \"573\u00d76=3438\";  After = \"872\u00d79=7848\" },\n    @{ Before = \"230\u00d72=460\";   After = \"423\u00d76=2538\" },\n    @{ Before = \"638\u00d78=5104\";  After = \"739\u00d75=3695\" },\n    @{ Before = \"216\u00d74=864\";   After = \"880\u00d76=5280\" },\n    @{ Before = \"599\u00d73=1797\";  After = \"183\u00d74=732\" },\n    @{ Before = \"824\u00d72=1648\";  After = \"294\u00d79=2646\" },\n    @{ Before = \"787\u00d78=6296\";  After = \"755\u00d76=4530\" },\n    @{ Before = \"815\u00d74=3260\";  After = \"837\u00d72=1674\" },\n    @{ Before = \"398\u00d72=796\";   After = \"663\u00d74=2652\" },\n    @{ Before = \"271\u00d77=1897\";  After = \"634\u00d78=5072\" },\n    @{ Before = \"943\u00d79=8487\";  After = \"611\u00d73=1833\" },\n    @{ Before = \"762\u00d76=4572\";  After = \"828\u00d74=3312\" },\n    @{ Before = \"406\u00d77=2842\";  After = \"958\u00d72=1916\" },\n    @{ Before = \"494\u00d78=3952\";  After = \"296\u00d74=1184\" },\n    @{ Before = \"924\u00d79=8316\";  After = \"165\u00d78=1320\" },\n    @{ Before = \"307\u00d79=2763\";  After = \"719\u00d75=3595\" },\n    @{ Before = \"239\u00d72=478\";   After = \"551\u00d79=4959\" },\n    @{ Before = \"846\u00d72=1692\";  After = \"664\u00d77=4648\" },\n    @{ Before = \"589\u00d72=1178\";  After = \"347\u00d78=2776\" },\n    @{ Before = \"740\u00d79=6660\";  After = \"649\u00d77=4543\" },\n    @{ Before = \"460\u00d74=1840\";  After = \"569\u00d75=2845\" },\n    @{ Before = \"439\u00d74=1756\";  After = \"240\u00d73=720\" },\n    @{ Before = \"391\u00d72=782\";   After = \"313\u00d72=626\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute(\n        $r.Before,\n        $true,    # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # wdFindContinue\n        $false,   # Format\n        $r.After,\n        2         # wdReplaceAll\n    )\n}\n"}
